$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Per-sending-cluster ligand stats: E,F,G,H,I,J
$ligandStats = @{
  "ECs" = @(3, 1, 1.305319666666667, 3.915959, 0.4586709810613518, 0.4586709810613519)
  "FAPs" = @(3, 1, 0.983774, 2.951322, 0.3456843539904149, 0.3456843539904149)
  "MuSCs" = @(1, 0.3333333333333333, 0.1949326666666667, 0.584798, 0.0684965987597716, 0.06849659875977161)
  "Resolving-Mac" = @(3, 1, 0.3618473333333334, 1.085542, 0.1271480661884616, 0.1271480661884616)
}

# Per-target-cluster receptor stats: K,L,M,N,O,P
$receptorStats = @{
  "ECs" = @(3, 1, 46.63275166666667, 139.898255, 0.9158911059585902, 0.9158911059585902)
  "FAPs" = @(3, 1, 2.770761666666667, 8.312285, 0.05441917700612491, 0.05441917700612491)
  "MuSCs" = @(3, 1, 0.849605, 2.548815, 0.01668667696558362, 0.01668667696558362)
  "Resolving-Mac" = @(3, 1, 0.662052, 1.986156, 0.01300304006970129, 0.0130030400697013)
}

# Per-(sending,target) pair edge stats: Q,R,S,T
$pairStats = @{
  "ECs|ECs" = @(60.87064786128278, 547.835830751545, 0.4200926721153931, 0.4200926721153932)
  "ECs|FAPs" = @(3.616729695146111, 32.550567256315, 0.02496049730595067, 0.02496049730595068)
  "ECs|MuSCs" = @(1.109006115398333, 9.981055038585001, 0.007653694494458103, 0.007653694494458104)
  "ECs|Resolving-Mac" = @(0.8641894959559999, 7.777705463604, 0.005964117145549962, 0.005964117145549962)
  "FAPs|ECs" = @(45.87608863812334, 412.88479774311, 0.3166092252888619, 0.3166092252888619)
  "FAPs|FAPs" = @(2.725803287863333, 24.53222959077, 0.01881185804805233, 0.01881185804805233)
  "FAPs|MuSCs" = @(0.83581930927, 7.522373783430002, 0.005768323147094512, 0.005768323147094512)
  "FAPs|Resolving-Mac" = @(0.651309544248, 5.861785898232, 0.00449494750640617, 0.004494947506406172)
  "MuSCs|ECs" = @(9.090246636387779, 81.81221972749, 0.06273542559248901, 0.06273542559248903)
  "MuSCs|FAPs" = @(0.5401119603811111, 4.86100764343, 0.003727528532225526, 0.003727528532225527)
  "MuSCs|MuSCs" = @(0.1656157682633334, 1.49054191437, 0.001142980616745505, 0.001142980616745505)
  "MuSCs|Resolving-Mac" = @(0.129055561832, 1.161500056488, 0.000890664018311562, 0.0008906640183115623)
  "Resolving-Mac|ECs" = @(16.87393683657889, 151.86543152921, 0.1164537829618462, 0.1164537829618462)
  "Resolving-Mac|FAPs" = @(1.002592720385556, 9.023334483469998, 0.006919293119896379, 0.00691929311989638)
  "Resolving-Mac|MuSCs" = @(0.3074273036366667, 2.76684573273, 0.002121678707285505, 0.002121678707285505)
  "Resolving-Mac|Resolving-Mac" = @(0.239561750728, 2.156055756552, 0.001653311399433599, 0.001653311399433599)
}

$clusters = @("ECs", "FAPs", "MuSCs", "Resolving-Mac")

$row = 2
foreach ($sending in $clusters) {
  foreach ($target in $clusters) {
    $ws.Cells.Item($row, 1).Value = $sending
    $ws.Cells.Item($row, 2).Value = "Dll3"
    $ws.Cells.Item($row, 3).Value = "Notch4"
    $ws.Cells.Item($row, 4).Value = $target

    $lig = $ligandStats[$sending]
    for ($i = 0; $i -lt 6; $i++) {
      $ws.Cells.Item($row, 5 + $i).Value = $lig[$i]
    }

    $rec = $receptorStats[$target]
    for ($i = 0; $i -lt 6; $i++) {
      $ws.Cells.Item($row, 11 + $i).Value = $rec[$i]
    }

    $pair = $pairStats["$sending|$target"]
    for ($i = 0; $i -lt 4; $i++) {
      $ws.Cells.Item($row, 17 + $i).Value = $pair[$i]
    }

    $row = $row + 1
  }
}